$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H64").Value = 4308.4736
$ws.Range("I64").Value = 3816.6667
$ws.Range("K64").Value = 3816.6667
$ws.Range("M64").Value = -3568.6667
$ws.Range("H67").Value = 4308.4736
$ws.Range("I67").Value = 3816.6667
$ws.Range("K67").Value = 3816.6667
$ws.Range("M67").Value = -2958.6667
$ws.Range("H115").Value = 200
$ws.Range("I115").Value = 200
$ws.Range("K115").Value = 600
$ws.Range("M115").Value = 967
$ws.Range("H129").Value = 871.8108
$ws.Range("J129").Value = 998.7
$ws.Range("L129").Value = 2996.1
$ws.Range("N129").Value = -12996.1
$ws.Range("H138").Value = 1170997.5
$ws.Range("J138").Value = 2290455.2
$ws.Range("L138").Value = 6871365.600000001
$ws.Range("N138").Value = -6881645.600000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23205.143
$ws.Range("I32").Value = 27116.238
$ws.Range("K32").Value = 27116.238
$ws.Range("M32").Value = -26829.238
$ws.Range("H35").Value = 2656.889
$ws.Range("I35").Value = 1558.7142
$ws.Range("J35").Value = 6500.5
$ws.Range("K35").Value = 1558.7142
$ws.Range("L35").Value = 6500.5
$ws.Range("M35").Value = -1152.7142
$ws.Range("N35").Value = -7312.5
$ws.Range("H43").Value = 8933
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 8933
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 8933
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -9559
$ws.Range("H76").Value = 38780
$ws.Range("J76").Value = 38780
$ws.Range("L76").Value = 38780
$ws.Range("N76").Value = -39456
$ws.Range("H79").Value = 38780
$ws.Range("J79").Value = 38780
$ws.Range("L79").Value = 38780
$ws.Range("N79").Value = -41120
$ws.Range("H106").Value = 60120.715
$ws.Range("J106").Value = 60120.715
$ws.Range("L106").Value = 60120.715
$ws.Range("N106").Value = -62644.715

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 39271.332
$ws.Range("J76").Value = 39271.332
$ws.Range("L76").Value = 39271.332
$ws.Range("N76").Value = -39901.332
$ws.Range("H79").Value = 39271.332
$ws.Range("J79").Value = 39271.332
$ws.Range("L79").Value = 39271.332
$ws.Range("N79").Value = -41455.332

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1485
$ws.Range("I16").Value = 1012.5
$ws.Range("J16").Value = 1957.5
$ws.Range("K16").Value = 1012.5
$ws.Range("L16").Value = 1957.5
$ws.Range("M16").Value = -725.5
$ws.Range("N16").Value = -2531.5
$ws.Range("H31").Value = 3221.5881
$ws.Range("I31").Value = 1737
$ws.Range("K31").Value = 1737
$ws.Range("M31").Value = -1442
$ws.Range("H34").Value = 3221.5881
$ws.Range("I34").Value = 1737
$ws.Range("K34").Value = 1737
$ws.Range("M34").Value = -1535
$ws.Range("H58").Value = 28573426
$ws.Range("J58").Value = 2745.875
$ws.Range("L58").Value = 2745.875
$ws.Range("N58").Value = -3151.875
$ws.Range("H112").Value = 32900.668
$ws.Range("J112").Value = 32900.668
$ws.Range("L112").Value = 32900.668
$ws.Range("N112").Value = -35854.668
$ws.Range("H113").Value = 1485
$ws.Range("I113").Value = 1012.5
$ws.Range("J113").Value = 1957.5
$ws.Range("K113").Value = 1012.5
$ws.Range("L113").Value = 1957.5
$ws.Range("M113").Value = 1157.5
$ws.Range("N113").Value = -6297.5
$ws.Range("H132").Value = 35526.367
$ws.Range("I132").Value = 2206.611
$ws.Range("K132").Value = 6619.833
$ws.Range("M132").Value = -4089.833
$ws.Range("H134").Value = 33409.734
$ws.Range("I134").Value = 1903.8695
$ws.Range("J134").Value = 99285.63
$ws.Range("K134").Value = 5711.6085
$ws.Range("L134").Value = 297856.89
$ws.Range("M134").Value = -3176.6085
$ws.Range("N134").Value = -302926.89
$ws.Range("H136").Value = 28573426
$ws.Range("J136").Value = 2745.875
$ws.Range("L136").Value = 8237.625
$ws.Range("N136").Value = -13337.625

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43.636364
$ws.Range("J2").Value = 22
$ws.Range("L2").Value = 132
$ws.Range("N2").Value = -358
$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 45000
$ws.Range("N42").Value = -46068
$ws.Range("H70").Value = 2751.8333
$ws.Range("I70").Value = 1007.3333
$ws.Range("J70").Value = 3333.3333
$ws.Range("K70").Value = 3021.9999
$ws.Range("L70").Value = 9999.999899999999
$ws.Range("M70").Value = -2706.9999
$ws.Range("N70").Value = -10629.9999
$ws.Range("H73").Value = 2751.8333
$ws.Range("I73").Value = 1007.3333
$ws.Range("J73").Value = 3333.3333
$ws.Range("K73").Value = 3021.9999
$ws.Range("L73").Value = 9999.999899999999
$ws.Range("M73").Value = -1929.9999
$ws.Range("N73").Value = -12183.9999
$ws.Range("H75").Value = 2198.818
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 2576.3333
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 7728.999899999999
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -9724.999899999999
$ws.Range("H78").Value = 2198.818
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 2576.3333
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 23186.9997
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -33170.9997
$ws.Range("H122").Value = 546.85
$ws.Range("I122").Value = 291.05884
$ws.Range("J122").Value = 1996.3334
$ws.Range("K122").Value = 2619.52956
$ws.Range("L122").Value = 17967.0006
$ws.Range("M122").Value = -169.5295599999999
$ws.Range("N122").Value = -22867.0006
$ws.Range("H132").Value = 1248.8096
$ws.Range("I132").Value = 908.9286
$ws.Range("K132").Value = 8180.3574
$ws.Range("M132").Value = -5650.3574

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 7666.2856
$ws.Range("I31").Value = 2277.3333
$ws.Range("K31").Value = 2277.3333
$ws.Range("M31").Value = -1985.3333
$ws.Range("H37").Value = 7666.2856
$ws.Range("I37").Value = 2277.3333
$ws.Range("K37").Value = 2277.3333
$ws.Range("M37").Value = -2000.3333
$ws.Range("H100").Value = 40750
$ws.Range("J100").Value = 40750
$ws.Range("L100").Value = 40750
$ws.Range("N100").Value = -42914

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 698.9375
$ws.Range("I9").Value = 192.11111
$ws.Range("J9").Value = 1350.5714
$ws.Range("K9").Value = 192.11111
$ws.Range("L9").Value = 1350.5714
$ws.Range("M9").Value = 31.88889
$ws.Range("N9").Value = -1798.5714
$ws.Range("H35").Value = 1443.6666
$ws.Range("I35").Value = 1443.6666
$ws.Range("K35").Value = 1443.6666
$ws.Range("M35").Value = -1107.6666
$ws.Range("H45").Value = 4000
$ws.Range("I45").Value = 4000
$ws.Range("K45").Value = 4000
$ws.Range("M45").Value = -3593
$ws.Range("H103").Value = 30643.75
$ws.Range("J103").Value = 30643.75
$ws.Range("L103").Value = 30643.75
$ws.Range("N103").Value = -32987.75
$ws.Range("H132").Value = 39826.777
$ws.Range("I132").Value = 1424.9286
$ws.Range("K132").Value = 4274.7858
$ws.Range("M132").Value = -1744.7858

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7250
$ws.Range("J41").Value = 7250
$ws.Range("L41").Value = 7250
$ws.Range("N41").Value = -8030
$ws.Range("H45").Value = 12999.667
$ws.Range("J45").Value = 12999.667
$ws.Range("L45").Value = 12999.667
$ws.Range("N45").Value = -13981.667
$ws.Range("H74").Value = 6626
$ws.Range("J74").Value = 6626
$ws.Range("L74").Value = 6626
$ws.Range("N74").Value = -8498
$ws.Range("H77").Value = 6626
$ws.Range("J77").Value = 6626
$ws.Range("L77").Value = 19878
$ws.Range("N77").Value = -29238
$ws.Range("H132").Value = 61794.91
$ws.Range("I132").Value = 36616.93
$ws.Range("K132").Value = 109850.79
$ws.Range("M132").Value = -107320.79
$ws.Range("H136").Value = 63696.156
$ws.Range("I136").Value = 38288.074
$ws.Range("K136").Value = 114864.222
$ws.Range("M136").Value = -112314.222
